$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("S2").Value = 1714

# Row 3
$ws.Range("D3").Value = 44320
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("S3").Value = 1143

# Row 4
$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 45

# Row 5
$ws.Range("D5").Value = 44322
$ws.Range("M5").Value = 80

# Row 8
$ws.Range("D8").Value = 44302
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 2143

# Row 9
$ws.Range("D9").Value = 44302
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("S9").Value = 1714

# Row 10
$ws.Range("D10").Value = 44300
$ws.Range("M10").Value = 100
$ws.Range("R10").Value = "Región Metropolitana"

# Row 11
$ws.Range("D11").Value = 44300
$ws.Range("M11").Value = 80
$ws.Range("R11").Value = "Región Metropolitana"

# Row 14
$ws.Range("D14").Value = 44299
$ws.Range("M14").Value = 80
$ws.Range("R14").Value = "Provincia de Santiago"

# Row 15
$ws.Range("D15").Value = 44299
$ws.Range("M15").Value = 75
$ws.Range("R15").Value = "Provincia de Santiago"
